$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep ID_Set (A) and nbPieces (C) columns as text, matching source inlineStr formatting
$ws.Range("A2:A26").NumberFormat = "@"
$ws.Range("C2:C26").NumberFormat = "@"

# Row 2: Set_10370
$ws.Range("A2").Value = '10370'
$ws.Range("B2").Value = 'L’étoile de Noël'
$ws.Range("C2").Value = '608'
$ws.Range("D2").Value = 'The Botanical Collection'
$ws.Range("E2").Value = 'https://www.lego.com/cdn/cs/set/assets/blt5b9064fcb12ba88b/10370_Prod.png?format=webply&fit=bounds&quality=75&width=528&height=528&dpr=1'
$ws.Range("F2").Value = 'https://www.lego.com/fr-fr/product/10370'
$ws.Range("G2").ClearContents()
$ws.Range("H2").ClearContents()
$ws.Range("I2").ClearContents()
$ws.Range("J2").ClearContents()

# Row 3: Set_10372
$ws.Range("A3").Value = '10372'
$ws.Range("B3").Value = 'Hibiscus'
$ws.Range("C3").Value = '660'
$ws.Range("D3").Value = 'The Botanical Collection'
$ws.Range("E3").Value = 'https://www.lego.com/cdn/cs/set/assets/blt4d1c7d7f1731540a/10372_Prod_en-gb.png?format=webply&fit=bounds&quality=75&width=800&height=800&dpr=1'
$ws.Range("F3").Value = 'https://www.lego.com/fr-fr/product/10372'
$ws.Range("G3").ClearContents()
$ws.Range("H3").ClearContents()
$ws.Range("I3").ClearContents()
$ws.Range("J3").ClearContents()

# Row 4: Set_11370
$ws.Range("A4").Value = '11370'
$ws.Range("B4").Value = 'Stranger Things : la Maison Creel'
$ws.Range("C4").Value = '2593'
$ws.Range("D4").Value = 'N/A'
$ws.Range("E4").Value = 'https://www.lego.com/cdn/cs/set/assets/blt8170cb2935ea6596/11370_Prod_en-gb.png?format=webply&fit=bounds&quality=75&width=528&height=528&dpr=1'
$ws.Range("F4").Value = 'https://www.lego.com/fr-fr/product/11370'
$ws.Range("G4").ClearContents()
$ws.Range("H4").ClearContents()
$ws.Range("I4").ClearContents()
$ws.Range("J4").ClearContents()

# Row 5: Set_11372
$ws.Range("A5").Value = '11372'
$ws.Range("B5").Value = 'Le jardin d''automne'
$ws.Range("C5").Value = '1102'
$ws.Range("D5").Value = 'N/A'
$ws.Range("E5").Value = 'https://www.lego.com/cdn/cs/set/assets/blt0345da5be8baeb2d/11372_Prod_en-gb.png?format=webply&fit=bounds&quality=75&width=528&height=528&dpr=1'
$ws.Range("F5").Value = 'https://www.lego.com/fr-fr/product/11372'
$ws.Range("G5").ClearContents()
$ws.Range("H5").ClearContents()
$ws.Range("I5").ClearContents()
$ws.Range("J5").ClearContents()

# Row 6: Set_11384
$ws.Range("A6").Value = '11384'
$ws.Range("B6").Value = 'Le chiot Golden Retriever'
$ws.Range("C6").Value = '2102'
$ws.Range("D6").Value = 'N/A'
$ws.Range("E6").Value = 'https://www.lego.com/cdn/cs/set/assets/bltc6a68358c38c3bf9/blt00e69bfd89abdb56-11384_Prod.png?format=webply&fit=bounds&quality=75&width=528&height=528&dpr=1'
$ws.Range("F6").Value = 'https://www.lego.com/fr-fr/product/11384'
$ws.Range("G6").ClearContents()
$ws.Range("H6").ClearContents()
$ws.Range("I6").ClearContents()
$ws.Range("J6").ClearContents()

# Row 7: Set_11508
$ws.Range("A7").Value = '11508'
$ws.Range("B7").Value = 'Marguerites'
$ws.Range("C7").Value = '133'
$ws.Range("D7").Value = 'N/A'
$ws.Range("E7").Value = 'https://www.lego.com/cdn/cs/set/assets/blt515cd2564f9b2dcf/11508_Prod_en-gb.png?format=webply&fit=bounds&quality=75&width=528&height=528&dpr=1'
$ws.Range("F7").Value = 'https://www.lego.com/fr-fr/product/11508'
$ws.Range("G7").ClearContents()
$ws.Range("H7").ClearContents()
$ws.Range("I7").ClearContents()
$ws.Range("J7").ClearContents()

# Row 8: Set_11509
$ws.Range("A8").Value = '11509'
$ws.Range("B8").Value = 'Cactus fleuri'
$ws.Range("C8").Value = '482'
$ws.Range("D8").Value = 'N/A'
$ws.Range("E8").Value = 'https://www.lego.com/cdn/cs/set/assets/blt05c278dbcfaa4871/11509_Prod_en-gb.png?format=webply&fit=bounds&quality=75&width=528&height=528&dpr=1'
$ws.Range("F8").Value = 'https://www.lego.com/fr-fr/product/11509'
$ws.Range("G8").ClearContents()
$ws.Range("H8").ClearContents()
$ws.Range("I8").ClearContents()
$ws.Range("J8").ClearContents()

# Row 9: Set_21365
$ws.Range("A9").Value = '21365'
$ws.Range("B9").Value = 'Les oiseaux amoureux'
$ws.Range("C9").Value = '750'
$ws.Range("D9").Value = 'N/A'
$ws.Range("E9").Value = 'https://www.lego.com/cdn/cs/set/assets/blt859e843a3f8c5af9/21365_Prod_en-gb.png?format=webply&fit=bounds&quality=75&width=528&height=528&dpr=1'
$ws.Range("F9").Value = 'https://www.lego.com/fr-fr/product/21365'
$ws.Range("G9").ClearContents()
$ws.Range("H9").ClearContents()
$ws.Range("I9").ClearContents()
$ws.Range("J9").ClearContents()

# Row 10: Set_31218
$ws.Range("A10").Value = '31218'
$ws.Range("B10").Value = 'Les cerisiers en fleurs'
$ws.Range("C10").Value = '1892'
$ws.Range("D10").Value = 'N/A'
$ws.Range("E10").Value = 'https://www.lego.com/cdn/cs/set/assets/blt3f6cace2ff7687c6/31218_Prod_en-gb.png?format=webply&fit=bounds&quality=75&width=528&height=528&dpr=1'
$ws.Range("F10").Value = 'https://www.lego.com/fr-fr/product/31218'
$ws.Range("G10").ClearContents()
$ws.Range("H10").ClearContents()
$ws.Range("I10").ClearContents()
$ws.Range("J10").ClearContents()

# Row 11: Set_31378
$ws.Range("A11").Value = '31378'
$ws.Range("B11").Value = 'Le télescope d’exploration spatiale'
$ws.Range("C11").Value = '278'
$ws.Range("D11").Value = 'N/A'
$ws.Range("E11").Value = 'https://www.lego.com/cdn/cs/set/assets/blt71fe323b45263a78/31378_Prod_en-gb.png?format=webply&fit=bounds&quality=75&width=528&height=528&dpr=1'
$ws.Range("F11").Value = 'https://www.lego.com/fr-fr/product/31378'
$ws.Range("G11").ClearContents()
$ws.Range("H11").ClearContents()
$ws.Range("I11").ClearContents()
$ws.Range("J11").ClearContents()

# Row 12: Set_31379
$ws.Range("A12").Value = '31379'
$ws.Range("B12").Value = 'Le dinosaure féroce'
$ws.Range("C12").Value = '283'
$ws.Range("D12").Value = 'N/A'
$ws.Range("E12").Value = 'https://www.lego.com/cdn/cs/set/assets/blt94ad0e9878a39d0a/31379_Prod_en-gb.png?format=webply&fit=bounds&quality=75&width=528&height=528&dpr=1'
$ws.Range("F12").Value = 'https://www.lego.com/fr-fr/product/31379'
$ws.Range("G12").ClearContents()
$ws.Range("H12").ClearContents()
$ws.Range("I12").ClearContents()
$ws.Range("J12").ClearContents()

# Row 13: Set_31384
$ws.Range("A13").Value = '31384'
$ws.Range("B13").Value = 'Animaux sauvages : le colibri coloré'
$ws.Range("C13").Value = '312'
$ws.Range("D13").Value = 'N/A'
$ws.Range("E13").Value = 'https://www.lego.com/cdn/cs/set/assets/blta735df411e869f69/31384_Prod_Crop.png?format=webply&fit=bounds&quality=75&width=528&height=528&dpr=1'
$ws.Range("F13").Value = 'https://www.lego.com/fr-fr/product/31384'
$ws.Range("G13").ClearContents()
$ws.Range("H13").ClearContents()
$ws.Range("I13").ClearContents()
$ws.Range("J13").ClearContents()

# Row 14: Set_31385
$ws.Range("A14").Value = '31385'
$ws.Range("B14").Value = 'Animaux marins : les beaux dauphins'
$ws.Range("C14").Value = '542'
$ws.Range("D14").Value = 'N/A'
$ws.Range("E14").Value = 'https://www.lego.com/cdn/cs/set/assets/bltbf6cbcbb8ff1610b/31385_Prod_en-gb.png?format=webply&fit=bounds&quality=75&width=528&height=528&dpr=1'
$ws.Range("F14").Value = 'https://www.lego.com/fr-fr/product/31385'
$ws.Range("G14").ClearContents()
$ws.Range("H14").ClearContents()
$ws.Range("I14").ClearContents()
$ws.Range("J14").ClearContents()

# Row 15: Set_40957
$ws.Range("A15").Value = '40957'
$ws.Range("B15").Value = 'La couronne printanière'
$ws.Range("C15").Value = '747'
$ws.Range("D15").Value = 'N/A'
$ws.Range("E15").Value = 'https://www.lego.com/cdn/cs/set/assets/bltd979e8936ef8f4e8/40957_Prod_en-gb.png?format=webply&fit=bounds&quality=75&width=528&height=528&dpr=1'
$ws.Range("F15").Value = 'https://www.lego.com/fr-fr/product/40957'
$ws.Range("G15").ClearContents()
$ws.Range("H15").ClearContents()
$ws.Range("I15").ClearContents()
$ws.Range("J15").ClearContents()

# Row 16: Set_42222
$ws.Range("A16").Value = '42222'
$ws.Range("B16").Value = 'Hypercar Bugatti Chiron Pur Sport'
$ws.Range("C16").Value = '771'
$ws.Range("D16").Value = 'N/A'
$ws.Range("E16").Value = 'https://www.lego.com/cdn/cs/set/assets/blt0a9009b95130ef89/42222_Prod_en-gb.png?format=webply&fit=bounds&quality=75&width=528&height=528&dpr=1'
$ws.Range("F16").Value = 'https://www.lego.com/fr-fr/product/42222'
$ws.Range("G16").ClearContents()
$ws.Range("H16").ClearContents()
$ws.Range("I16").ClearContents()
$ws.Range("J16").ClearContents()

# Row 17: Set_42227
$ws.Range("A17").Value = '42227'
$ws.Range("B17").Value = 'SUV Jeep® Wrangler Rubicon'
$ws.Range("C17").Value = '723'
$ws.Range("D17").Value = 'N/A'
$ws.Range("E17").Value = 'https://www.lego.com/cdn/cs/set/assets/blt8f7d3695517d3132/42227_Prod.png?format=webply&fit=bounds&quality=75&width=528&height=528&dpr=1'
$ws.Range("F17").Value = 'https://www.lego.com/fr-fr/product/42227'
$ws.Range("G17").ClearContents()
$ws.Range("H17").ClearContents()
$ws.Range("I17").ClearContents()
$ws.Range("J17").ClearContents()

# Row 18: Set_42684
$ws.Range("A18").Value = '42684'
$ws.Range("B18").Value = 'Le café licorne'
$ws.Range("C18").Value = '475'
$ws.Range("D18").Value = 'N/A'
$ws.Range("E18").Value = 'https://www.lego.com/cdn/cs/set/assets/blt83ef724ade6083c3/42684_Prod_en-gb.png?format=webply&fit=bounds&quality=75&width=528&height=528&dpr=1'
$ws.Range("F18").Value = 'https://www.lego.com/fr-fr/product/42684'
$ws.Range("G18").ClearContents()
$ws.Range("H18").ClearContents()
$ws.Range("I18").ClearContents()
$ws.Range("J18").ClearContents()

# Row 19: Set_42688
$ws.Range("A19").Value = '42688'
$ws.Range("B19").Value = 'L’écurie et l’école d’équitation'
$ws.Range("C19").Value = '735'
$ws.Range("D19").Value = 'N/A'
$ws.Range("E19").Value = 'https://www.lego.com/cdn/cs/set/assets/blt681e22a3e533a70a/42688_Prod_en-gb.png?format=webply&fit=bounds&quality=75&width=528&height=528&dpr=1'
$ws.Range("F19").Value = 'https://www.lego.com/fr-fr/product/42688'
$ws.Range("G19").ClearContents()
$ws.Range("H19").ClearContents()
$ws.Range("I19").ClearContents()
$ws.Range("J19").ClearContents()

# Row 20: Set_42696
$ws.Range("A20").Value = '42696'
$ws.Range("B20").Value = 'La clinique vétérinaire'
$ws.Range("C20").Value = '141'
$ws.Range("D20").Value = 'N/A'
$ws.Range("E20").Value = 'https://www.lego.com/cdn/cs/set/assets/bltc4781be88c8bb640/42696_Prod_en-gb.png?format=webply&fit=bounds&quality=75&width=528&height=528&dpr=1'
$ws.Range("F20").Value = 'https://www.lego.com/fr-fr/product/42696'
$ws.Range("G20").ClearContents()
$ws.Range("H20").ClearContents()
$ws.Range("I20").ClearContents()
$ws.Range("J20").ClearContents()

# Row 21: Set_43281
$ws.Range("A21").Value = '43281'
$ws.Range("B21").Value = 'Le château de glace et la piste enneigée d’Elsa'
$ws.Range("C21").Value = '216'
$ws.Range("D21").Value = 'N/A'
$ws.Range("E21").Value = 'https://www.lego.com/cdn/cs/set/assets/blt3c56a41db2717b7b/43281_Prod_en-gb.png?format=webply&fit=bounds&quality=75&width=528&height=528&dpr=1'
$ws.Range("F21").Value = 'https://www.lego.com/fr-fr/product/43281'
$ws.Range("G21").ClearContents()
$ws.Range("H21").ClearContents()
$ws.Range("I21").ClearContents()
$ws.Range("J21").ClearContents()

# Row 22: Set_43287
$ws.Range("A22").Value = '43287'
$ws.Range("B22").Value = 'Le pique-nique d’Olaf et Bruni'
$ws.Range("C22").Value = '478'
$ws.Range("D22").Value = 'N/A'
$ws.Range("E22").Value = 'https://www.lego.com/cdn/cs/set/assets/blt9bcf160a34a7c824/43287_Prod_en-gb.png?format=webply&fit=bounds&quality=75&width=528&height=528&dpr=1'
$ws.Range("F22").Value = 'https://www.lego.com/fr-fr/product/43287'
$ws.Range("G22").ClearContents()
$ws.Range("H22").ClearContents()
$ws.Range("I22").ClearContents()
$ws.Range("J22").ClearContents()

# Row 23: Set_45200
$ws.Range("A23").Value = '45200'
$ws.Range("B23").Value = 'Kit de science Mission lunaire'
$ws.Range("C23").Value = '519'
$ws.Range("D23").Value = 'N/A'
$ws.Range("E23").Value = 'https://www.lego.com/cdn/cs/set/assets/bltf7b7ce8f03c1660e/bltb54bf9af8b6f95a9-45200_Prod.png?format=webply&fit=bounds&quality=75&width=528&height=528&dpr=1'
$ws.Range("F23").Value = 'https://www.lego.com/fr-fr/product/45200'
$ws.Range("G23").ClearContents()
$ws.Range("H23").ClearContents()
$ws.Range("I23").ClearContents()
$ws.Range("J23").ClearContents()

# Row 24: Set_75423
$ws.Range("A24").Value = '75423'
$ws.Range("B24").Value = 'SMART Play™ : le X-Wing™ Red Five de Luke'
$ws.Range("C24").Value = '581'
$ws.Range("D24").Value = 'N/A'
$ws.Range("E24").Value = 'https://www.lego.com/cdn/cs/set/assets/blt5083a28189d02095/75423_Prod.png?format=webply&fit=bounds&quality=75&width=528&height=528&dpr=1'
$ws.Range("F24").Value = 'https://www.lego.com/fr-fr/product/75423'
$ws.Range("G24").ClearContents()
$ws.Range("H24").ClearContents()
$ws.Range("I24").ClearContents()
$ws.Range("J24").ClearContents()

# Row 25: Set_75440
$ws.Range("A25").Value = '75440'
$ws.Range("B25").Value = 'AT-AT™'
$ws.Range("C25").Value = '525'
$ws.Range("D25").Value = 'N/A'
$ws.Range("E25").Value = 'https://www.lego.com/cdn/cs/set/assets/blt26f813f1a1f499c6/75440_Prod_en-gb.png?format=webply&fit=bounds&quality=75&width=528&height=528&dpr=1'
$ws.Range("F25").Value = 'https://www.lego.com/fr-fr/product/75440'
$ws.Range("G25").ClearContents()
$ws.Range("H25").ClearContents()
$ws.Range("I25").ClearContents()
$ws.Range("J25").ClearContents()

# Row 26: Set_77256
$ws.Range("A26").Value = '77256'
$ws.Range("B26").Value = 'La machine à remonter le temps de Retour vers le futur'
$ws.Range("C26").Value = '357'
$ws.Range("D26").Value = 'N/A'
$ws.Range("E26").Value = 'https://www.lego.com/cdn/cs/set/assets/bltcf1faba2d5f4a688/77256_Prod_en-gb.png?format=webply&fit=bounds&quality=75&width=528&height=528&dpr=1'
$ws.Range("F26").Value = 'https://www.lego.com/fr-fr/product/77256'
$ws.Range("G26").ClearContents()
$ws.Range("H26").ClearContents()
$ws.Range("I26").ClearContents()
$ws.Range("J26").ClearContents()

# Remove now-unused trailing rows (27-29) -- dimension shrinks to A1:J26
$ws.Range("A27:J29").Delete()
